$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 (Instar, Month 6)
$ws.Range("C4").Value = 5.97371428571429
$ws.Range("D4").Value = 175

# Row 5 (Instar, Month 7)
$ws.Range("C5").Value = 5.19856
$ws.Range("D5").Value = 125

# Row 9 (Megalopae, Month 6)
$ws.Range("C9").Value = 2.77320430107527
$ws.Range("D9").Value = 465

# Row 10 (Megalopae, Month 7)
$ws.Range("C10").Value = 2.42452153110048
$ws.Range("D10").Value = 419

# Row 11 (Megalopae, Month 8)
$ws.Range("C11").Value = 2.54321428571429
$ws.Range("D11").Value = 57
